# Apply the FHIR IG StructureDefinition metadata refresh:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date refresh
#  - Publisher set to "Alvearie Team"
#  - Contact (x2 rows) replaced by a single Jurisdiction row
#  - Stray "N/A" note on the parameterValue extension row cleared

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Replace the first "Contact" row (row 10) with "Jurisdiction"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Remove the now-duplicate second "Contact" row (row 11), shifting rows 12-20 up
$meta.Rows.Item(11).Delete()

# Clear the stray "N/A" mapping note on the parameterValue extension row
$elements.Cells.Item(5, 36).Value = ""
